$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 'Human Factors in Cyber Warfare II: Emerging Perspectives'
$ws.Range("C2").Value = 'Panel Chair: Dr.Vincent F. Mancuso, Panelists: Dr.James C. Christensen, Dr.Jennifer Cowley, Dr.Victor Finomore, Prof.Cleotide Gonzalez, Dr.Benjamin Knott'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '2014'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '10.1177/1541931214581085'
$ws.Range("F2").Value = 'Restricted'

# Row 3
$ws.Range("B3").Value = 'Modeling and application for pneumatic soft actuators based on a novel deep neural network'
$ws.Range("C3").Value = 'Ke Zhang, Yongqi Bi, Ruiyu Zhang'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2025'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '10.1177/09544062251317813'
$ws.Range("F3").Value = 'Restricted'

# Row 4
$ws.Range("B4").Value = 'Technical note: exploiting problem definition study for cyber security simulations'
$ws.Range("C4").Value = 'Yilmaz Cankaya'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2015'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '10.1177/1548512915604585'

# Row 5
$ws.Range("B5").Value = 'Cyber-Flirting: Playing at Love on the Internet'
$ws.Range("C5").Value = 'Monica Therese Whitty'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '2003'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '10.1177/0959354303013003003'
$ws.Range("G5").Value = 3

# Row 6
$ws.Range("B6").Value = 'Prioritizing investment in military cyber capability using risk analysis'
$ws.Range("C6").Value = 'Cayt Rowe, Hossein Seif Zadeh, Ivan L. Garanovich, Li Jiang, Daniel Bilusich, Rick Nunes-Vaz, Anthony Ween'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2019'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '10.1177/1548512917707077'

# Row 7
$ws.Range("B7").Value = 'How the process of discovering cyberattacks biases our understanding of cybersecurity'
$ws.Range("C7").Value = 'Harry Oppenheimer'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2024'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '10.1177/00223433231217687'
$ws.Range("F7").Value = 'Open Access'

# Row 8
$ws.Range("B8").Value = 'Ontological security, cyber technology, and states’ responses'
$ws.Range("C8").Value = 'Amir Lupovici'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2023'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '10.1177/13540661221130958'
$ws.Range("F8").Value = 'Open Access'
$ws.Range("G8").Value = 0

# Row 9
$ws.Range("B9").Value = 'Evaluation of communIT, a large-scale, cyber-physical artifact supporting diverse subgroups building community'
$ws.Range("C9").Value = 'Carlos de Aguiar, Gilly Leshed, Trevor Pinch, Keith Green'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2022'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '10.3233/SCS-220007'

# Row 10
$ws.Range("B10").Value = 'Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators'
$ws.Range("C10").Value = 'Kazem Kazerounian, Zhaoyu Wang'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1988'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '10.1177/027836498800700501'
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("B11").Value = 'The dynamics of cyber conflict between rival antagonists, 2001–11'
$ws.Range("C11").Value = 'Brandon Valeriano, Ryan C Maness'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2014'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10.1177/0022343313518940'
$ws.Range("G11").Value = 1
